# Update cryptos list with latest prices / 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.187.20"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "'3.044.21"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'595.80"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'154.90"
$ws.Range("E6").Value = "  +7.74%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'3.046.52"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "'6.91"
$ws.Range("E10").Value = "  +14.23%  "
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "'0.0000236"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "'35.93"
$ws.Range("E14").Value = "  +4.44%  "
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "'3.545.46"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "'7.12"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "'63.117.36"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "'3.044.44"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "'456.55"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'14.35"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'83.16"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'11.26"
$ws.Range("E25").Value = "  +4.66%  "
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").Value = "'12.39"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +4.33%  "
$ws.Range("D30").Value = "'2.25"
$ws.Range("E30").Value = "  +9.11%  "
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'27.78"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("D37").Value = "'5.96"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").Value = "'3.21"
$ws.Range("E38").Value = "  +11.62%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.132"
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.13"
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("D41").Value = "'50.47"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'9.16"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "'0.304"
$ws.Range("E43").Value = "  +12.21%  "
$ws.Range("D44").Value = "'43.77"
$ws.Range("E44").Value = "  +9.53%  "
$ws.Range("D45").Value = "'394.87"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").Value = "'2.730.15"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "'132.27"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +6.60%  "
$ws.Range("D51").Value = "'24.60"
$ws.Range("E51").Value = "  +4.42%  "
